$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(1, 1).Value = "Organisatie "
$ws.Cells.Item(1, 2).Value = "Status inleveren routekaart"

$ws.Cells.Item(2, 1).Value = "Aafje Thuiszorg Huizen Zorghotels (Stichting)"
$ws.Cells.Item(2, 2).Value = "Actueel en vastgesteld"
$ws.Cells.Item(3, 1).Value = "Accolade (Stichting)"
$ws.Cells.Item(3, 2).Value = "Wordt aan gewerkt"
$ws.Cells.Item(4, 1).Value = "ActiVite (Stichting)"
$ws.Cells.Item(4, 2).Value = "Wordt aan gewerkt"
$ws.Cells.Item(5, 1).Value = "Adullam Zorg (Stichting)"
$ws.Cells.Item(5, 2).Value = "Actueel en vastgesteld"
$ws.Cells.Item(6, 1).Value = "Alkcare (Stichting)"
$ws.Cells.Item(6, 2).Value = "Actueel en vastgesteld"
$ws.Cells.Item(7, 1).Value = "Alliade "
$ws.Cells.Item(7, 2).Value = "Actueel en vastgesteld"
$ws.Cells.Item(8, 1).Value = "Altrecht (Stichting)"
$ws.Cells.Item(8, 2).Value = "Actueel en vastgesteld"
$ws.Cells.Item(9, 1).Value = "Amaliazorg (Stichting)"
$ws.Cells.Item(9, 2).Value = "Actueel en vastgesteld"
$ws.Cells.Item(10, 1).Value = "Amarant Groep (Stichting)"
$ws.Cells.Item(10, 2).Value = "Actueel en vastgesteld"
$ws.Cells.Item(11, 1).Value = "Amaris Zorggroep (Stichting)"
$ws.Cells.Item(11, 2).Value = "Wordt aan gewerkt"
$ws.Cells.Item(12, 1).Value = "Amerpoort (Stichting)"
$ws.Cells.Item(12, 2).Value = "Wordt aan gewerkt"
$ws.Cells.Item(13, 1).Value = "Amstelring Groep (Stichting)"
$ws.Cells.Item(13, 2).Value = "Wordt aan gewerkt"
$ws.Cells.Item(14, 1).Value = "Archipel (Stichting)"
$ws.Cells.Item(14, 2).Value = "Actueel en vastgesteld"
$ws.Cells.Item(15, 1).Value = "Argos Zorggroep (Stichting)"
$ws.Cells.Item(15, 2).Value = "Wordt aan gewerkt"
$ws.Cells.Item(16, 1).Value = "Arkin (Stichting)"
$ws.Cells.Item(16, 2).Value = "Actueel en vastgesteld"
$ws.Cells.Item(17, 1).Value = "ASVZ (Stichting)"
$ws.Cells.Item(17, 2).Value = "Wordt aan gewerkt"
$ws.Cells.Item(18, 1).Value = "Atlant Zorggroep (Stichting)"
$ws.Cells.Item(18, 2).Value = "Wordt aan gewerkt"
$ws.Cells.Item(19, 1).Value = "Attent"
$ws.Cells.Item(19, 2).Value = "Wordt aan gewerkt"
$ws.Cells.Item(20, 1).Value = "Aveleijn (Stichting)"
$ws.Cells.Item(20, 2).Value = "Wordt aan gewerkt"
$ws.Cells.Item(21, 1).Value = "AxionContinu Groep (Stichting)"
$ws.Cells.Item(21, 2).Value = "Actueel en vastgesteld"
$ws.Cells.Item(22, 1).Value = "Azora (Stichting)"
$ws.Cells.Item(22, 2).Value = "Wordt aan gewerkt"
$ws.Cells.Item(23, 1).Value = "Baalderborg Groep (Stichting)"
$ws.Cells.Item(23, 2).Value = "Wordt aan gewerkt"
$ws.Cells.Item(24, 1).Value = "Bartholomeus Gasthuis (Stichting)"
$ws.Cells.Item(24, 2).Value = "Wordt aan gewerkt"
$ws.Cells.Item(25, 1).Value = "Bartimeus"
$ws.Cells.Item(25, 2).Value = "Actueel en vastgesteld"
$ws.Cells.Item(26, 1).Value = "Bethanië (Stichting)"
$ws.Cells.Item(26, 2).Value = "Actueel en vastgesteld"
$ws.Cells.Item(27, 1).Value = "Beweging 3.0 (Stichting)"
$ws.Cells.Item(27, 2).Value = "Actueel en vastgesteld"
$ws.Cells.Item(28, 1).Value = "BrabantZorg"
$ws.Cells.Item(28, 2).Value = "Actueel en vastgesteld"
$ws.Cells.Item(29, 1).Value = "Breederzorg Thuiszorg (Stichting)"
$ws.Cells.Item(29, 2).Value = "Actueel en vastgesteld"
$ws.Cells.Item(30, 1).Value = "Cardia (Stichting)"
$ws.Cells.Item(30, 2).Value = "Actueel en vastgesteld"
$ws.Cells.Item(31, 1).Value = "Careander (Stichting)"
$ws.Cells.Item(31, 2).Value = "Wordt aan gewerkt"
$ws.Cells.Item(32, 1).Value = "Careyn (Stichting)"
$ws.Cells.Item(32, 2).Value = "Actueel en vastgesteld"
$ws.Cells.Item(33, 1).Value = "Carinova (Stichting)"
$ws.Cells.Item(33, 2).Value = "Wordt aan gewerkt"
$ws.Cells.Item(34, 1).Value = "Carint-Reggeland Groep (Stichting)"
$ws.Cells.Item(34, 2).Value = "Actueel en vastgesteld"
$ws.Cells.Item(35, 1).Value = "Cedrah (Stichting)"
$ws.Cells.Item(35, 2).Value = "Actueel en vastgesteld"
$ws.Cells.Item(36, 1).Value = "Charim (Zorggroep)"
$ws.Cells.Item(36, 2).Value = "Wordt aan gewerkt"
$ws.Cells.Item(37, 1).Value = "Cicero Zorggroep (Stichting)"
$ws.Cells.Item(37, 2).Value = "Actueel en vastgesteld"
$ws.Cells.Item(38, 1).Value = "Combinatie Jeugdzorg"
$ws.Cells.Item(38, 2).Value = "Wordt aan gewerkt"
$ws.Cells.Item(39, 1).Value = "Cordaan Groep (Stichting)"
$ws.Cells.Item(39, 2).Value = "Wordt aan gewerkt"
$ws.Cells.Item(40, 1).Value = "Cosis"
$ws.Cells.Item(40, 2).Value = "Wordt aan gewerkt"
$ws.Cells.Item(41, 1).Value = "Curamare (Stichting)"
$ws.Cells.Item(41, 2).Value = "Wordt aan gewerkt"
$ws.Cells.Item(42, 1).Value = "Dagelijks Leven Zorg (B.V)."
$ws.Cells.Item(42, 2).Value = "Actueel en vastgesteld"
$ws.Cells.Item(43, 1).Value = "De Blije Borgh / Protestants Interkerkelijke Stichting Zorgverlening Hendrik-Ido-Ambacht"
$ws.Cells.Item(43, 2).Value = "Actueel en vastgesteld"
$ws.Cells.Item(44, 1).Value = "De Hoop ggz (Stichting)"
$ws.Cells.Item(44, 2).Value = "Actueel en vastgesteld"
$ws.Cells.Item(45, 1).Value = "De Hoven (Stichting)"
$ws.Cells.Item(45, 2).Value = "Wordt aan gewerkt"
$ws.Cells.Item(46, 1).Value = "De Lange Wei (Stichting) (Burgemeester De Boer)"
$ws.Cells.Item(46, 2).Value = "Actueel en vastgesteld"
$ws.Cells.Item(47, 1).Value = "De Wijngaerd (Stichting)"
$ws.Cells.Item(47, 2).Value = "Actueel en vastgesteld"
$ws.Cells.Item(48, 1).Value = "De Zijlen (Stichting) (Ilmarinen)"
$ws.Cells.Item(48, 2).Value = "Actueel en vastgesteld"
$ws.Cells.Item(49, 1).Value = "De Zorgcirkel (Stichting)"
$ws.Cells.Item(49, 2).Value = "Actueel en vastgesteld"
$ws.Cells.Item(50, 1).Value = "DFZS De Forensische Zorgspecialisten"
$ws.Cells.Item(50, 2).Value = "Actueel en vastgesteld"
$ws.Cells.Item(51, 1).Value = "Dichterbij (Stichting)"
$ws.Cells.Item(51, 2).Value = "Actueel en vastgesteld"
$ws.Cells.Item(52, 1).Value = "Dimence Groep (Stichting)"
$ws.Cells.Item(52, 2).Value = "Actueel en vastgesteld"
$ws.Cells.Item(53, 1).Value = "Domus Magnus B.V."
$ws.Cells.Item(53, 2).Value = "Actueel en vastgesteld"
$ws.Cells.Item(54, 1).Value = "Driegasthuizengroep"
$ws.Cells.Item(54, 2).Value = "Wordt aan gewerkt"
$ws.Cells.Item(55, 1).Value = "Driestroom (Stichting)"
$ws.Cells.Item(55, 2).Value = "Wordt aan gewerkt"
$ws.Cells.Item(56, 1).Value = "DSV (Stichting)"
$ws.Cells.Item(56, 2).Value = "Actueel en vastgesteld"
$ws.Cells.Item(57, 1).Value = "Elver (Stichting)"
$ws.Cells.Item(57, 2).Value = "Wordt aan gewerkt"
$ws.Cells.Item(58, 1).Value = "Espria (stichting ...) (met onderdelen Trans, Meander, Evean, GGZ Drenthe, icare)"
$ws.Cells.Item(58, 2).Value = "Actueel en vastgesteld"
$ws.Cells.Item(59, 1).Value = "Fier Fryslan (Stichting)"
$ws.Cells.Item(59, 2).Value = "Actueel en vastgesteld"
$ws.Cells.Item(60, 1).Value = "Fivoor (forensische en intensieve psychiatrie) (Z-H)"
$ws.Cells.Item(60, 2).Value = "Actueel en vastgesteld"
$ws.Cells.Item(61, 1).Value = "Fundis (Stichting) (beheren meerdere zorginstellingen, zoals Welthuis)"
$ws.Cells.Item(61, 2).Value = "Actueel en vastgesteld"
$ws.Cells.Item(62, 1).Value = "GGNet (Stichting)"
$ws.Cells.Item(62, 2).Value = "Wordt aan gewerkt"
$ws.Cells.Item(63, 1).Value = "GGz Breburg Groep (Stichting)"
$ws.Cells.Item(63, 2).Value = "Actueel en vastgesteld"
$ws.Cells.Item(64, 1).Value = "GGZ Delfland - Geestelijke Gezondheidszorg Delfland (Stichting)"
$ws.Cells.Item(64, 2).Value = "Actueel en vastgesteld"
$ws.Cells.Item(65, 1).Value = "GGZ Friesland"
$ws.Cells.Item(65, 2).Value = "Actueel en vastgesteld"
$ws.Cells.Item(66, 1).Value = "GGZ inGeest (Stichting)"
$ws.Cells.Item(66, 2).Value = "Actueel en vastgesteld"
$ws.Cells.Item(67, 1).Value = "GGzE (Stichting) GGz Eindhoven"
$ws.Cells.Item(67, 2).Value = "Wordt aan gewerkt"
$ws.Cells.Item(68, 1).Value = "Het Gastenhuis B.V."
$ws.Cells.Item(68, 2).Value = "Actueel en vastgesteld"
$ws.Cells.Item(69, 1).Value = "Het Gasthuis Millingen aan de Rijn (Sint Jan De Deo) (Stichting)"
$ws.Cells.Item(69, 2).Value = "Actueel en vastgesteld"
$ws.Cells.Item(70, 1).Value = "Het Laar (Stichting)"
$ws.Cells.Item(70, 2).Value = "Actueel en vastgesteld"
$ws.Cells.Item(71, 1).Value = "Het Parkhuis (Stichting)"
$ws.Cells.Item(71, 2).Value = "Wordt aan gewerkt"
$ws.Cells.Item(72, 1).Value = "Hillegomse Organisaties voor Zorgverlening aan Ouderen / HOZO (Stichting)"
$ws.Cells.Item(72, 2).Value = "Actueel en vastgesteld"
$ws.Cells.Item(73, 1).Value = "HilverZorg (Stichting)"
$ws.Cells.Item(73, 2).Value = "Wordt aan gewerkt"
$ws.Cells.Item(74, 1).Value = "Huis Ter Leede (Protestantse Interkerkelijke Stichting)"
$ws.Cells.Item(74, 2).Value = "Actueel en vastgesteld"
$ws.Cells.Item(75, 1).Value = "Huize Winterdijk (Stichting tot Oprichting en Instandhouding van Bejaardenoorden en Verzorgingstehuizen, uitgaande v.d. Gereformeerde Gemeente in Nederland)"
$ws.Cells.Item(75, 2).Value = "Actueel en vastgesteld"
$ws.Cells.Item(76, 1).Value = "Humanitas (Stichting)"
$ws.Cells.Item(76, 2).Value = "Wordt aan gewerkt"
$ws.Cells.Item(77, 1).Value = "Humanitas voor Dienstverlening aan Mensen met een Handicap / Humanitas DMH (Stichting)"
$ws.Cells.Item(77, 2).Value = "Wordt aan gewerkt"
$ws.Cells.Item(78, 1).Value = "Innoforte (Stichting)"
$ws.Cells.Item(78, 2).Value = "Wordt aan gewerkt"
$ws.Cells.Item(79, 1).Value = "Interzorg Noord Nederland (Stichting)"
$ws.Cells.Item(79, 2).Value = "Actueel en vastgesteld"
$ws.Cells.Item(80, 1).Value = "Ipse de Bruggen (Stichting)"
$ws.Cells.Item(80, 2).Value = "Actueel en vastgesteld"
$ws.Cells.Item(81, 1).Value = "IrisZorg"
$ws.Cells.Item(81, 2).Value = "Actueel en vastgesteld"
$ws.Cells.Item(82, 1).Value = "Joris Zorg (Stichting)"
$ws.Cells.Item(82, 2).Value = "Wordt aan gewerkt"
$ws.Cells.Item(83, 1).Value = "Kalorama (Stichting)"
$ws.Cells.Item(83, 2).Value = "Actueel en vastgesteld"
$ws.Cells.Item(84, 1).Value = "Karakter (Stichting)"
$ws.Cells.Item(84, 2).Value = "Actueel en vastgesteld"
$ws.Cells.Item(85, 1).Value = "Kempenhaeghe (Stichting)"
$ws.Cells.Item(85, 2).Value = "Wordt aan gewerkt"
$ws.Cells.Item(86, 1).Value = "Kennemerhart"
$ws.Cells.Item(86, 2).Value = "Actueel en vastgesteld"
$ws.Cells.Item(87, 1).Value = "Kenter Jeugdhulp"
$ws.Cells.Item(87, 2).Value = "Wordt aan gewerkt"
$ws.Cells.Item(88, 1).Value = "Klein Geluk"
$ws.Cells.Item(88, 2).Value = "Actueel en vastgesteld"
$ws.Cells.Item(89, 1).Value = "Koninklijke Kentalis (Stichting)"
$ws.Cells.Item(89, 2).Value = "Actueel en vastgesteld"
$ws.Cells.Item(90, 1).Value = "Koninklijke Visio, expertisecentrum voor slechtziende en blinde mensen (Stichting)"
$ws.Cells.Item(90, 2).Value = "Wordt aan gewerkt"
$ws.Cells.Item(91, 1).Value = "Koperhorst (Stichting)"
$ws.Cells.Item(91, 2).Value = "Wordt aan gewerkt"
$ws.Cells.Item(92, 1).Value = "Koraal Groep (Stichting)"
$ws.Cells.Item(92, 2).Value = "Wordt aan gewerkt"
$ws.Cells.Item(93, 1).Value = "Korian Zorg B.V."
$ws.Cells.Item(93, 2).Value = "Actueel en vastgesteld"
$ws.Cells.Item(94, 1).Value = "KwadrantGroep (Stichting)"
$ws.Cells.Item(94, 2).Value = "Actueel en vastgesteld"
$ws.Cells.Item(95, 1).Value = "Land van Horne (Stichting voor Verpleeg-, Verzorgings- en Woonfaciliteiten ...)"
$ws.Cells.Item(95, 2).Value = "Actueel en vastgesteld"
$ws.Cells.Item(96, 1).Value = "Landelijke Stichting Vredenoord"
$ws.Cells.Item(96, 2).Value = "Actueel en vastgesteld"
$ws.Cells.Item(97, 1).Value = "Laurens (Stichting)"
$ws.Cells.Item(97, 2).Value = "Actueel en vastgesteld"
$ws.Cells.Item(98, 1).Value = "Leger des Heils Welzijns- en Gezondheidszorg (Stichting)"
$ws.Cells.Item(98, 2).Value = "Actueel en vastgesteld"
$ws.Cells.Item(99, 1).Value = "Lelie Zorggroep (Stichting)"
$ws.Cells.Item(99, 2).Value = "Actueel en vastgesteld"
$ws.Cells.Item(100, 1).Value = "Lentekind (Stichting)"
$ws.Cells.Item(100, 2).Value = "Actueel en vastgesteld"
$ws.Cells.Item(101, 1).Value = "Lentis incl. Dignis"
$ws.Cells.Item(101, 2).Value = "Wordt aan gewerkt"
$ws.Cells.Item(102, 1).Value = "Levvel (noord holland)"
$ws.Cells.Item(102, 2).Value = "Wordt aan gewerkt"
$ws.Cells.Item(103, 1).Value = "Liante (Stichting)"
$ws.Cells.Item(103, 2).Value = "Wordt aan gewerkt"
$ws.Cells.Item(104, 1).Value = "Liemerije (Stichting)"
$ws.Cells.Item(104, 2).Value = "Actueel en vastgesteld"
$ws.Cells.Item(105, 1).Value = "Livio (Stichting)"
$ws.Cells.Item(105, 2).Value = "Actueel en vastgesteld"
$ws.Cells.Item(106, 1).Value = "Magenta"
$ws.Cells.Item(106, 2).Value = "Wordt aan gewerkt"
$ws.Cells.Item(107, 1).Value = "Marente (Stichting)"
$ws.Cells.Item(107, 2).Value = "Wordt aan gewerkt"
$ws.Cells.Item(108, 1).Value = "Mariënstede (Stichting) (incl. Vughterstede)"
$ws.Cells.Item(108, 2).Value = "Actueel en vastgesteld"
$ws.Cells.Item(109, 1).Value = "Mediant, Stichting voor Geestelijke Gezondheidszorg Oost- en Midden Twente"
$ws.Cells.Item(109, 2).Value = "Wordt aan gewerkt"
$ws.Cells.Item(110, 1).Value = "MET-GGZ (Limburg)"
$ws.Cells.Item(110, 2).Value = "Actueel en vastgesteld"
$ws.Cells.Item(111, 1).Value = "Middin (Stichting)"
$ws.Cells.Item(111, 2).Value = "Wordt aan gewerkt"
$ws.Cells.Item(112, 1).Value = "Mondriaan (Stichting)"
$ws.Cells.Item(112, 2).Value = "Actueel en vastgesteld"
$ws.Cells.Item(113, 1).Value = "Nieuw Woelwijck, Dorpsgemeenschap van Geestelijk Gehandicapten (Stichting)"
$ws.Cells.Item(113, 2).Value = "Wordt aan gewerkt"
$ws.Cells.Item(114, 1).Value = "NNCZ (Noord Nederlandse Coöperatie van Zorgorganisaties)"
$ws.Cells.Item(114, 2).Value = "Actueel en vastgesteld"
$ws.Cells.Item(115, 1).Value = "Noorderboog (Stichting)"
$ws.Cells.Item(115, 2).Value = "Wordt aan gewerkt"
$ws.Cells.Item(116, 1).Value = "Noorderbreedte"
$ws.Cells.Item(116, 2).Value = "Actueel en vastgesteld"
$ws.Cells.Item(117, 1).Value = "Novadic-Kentron (Stichting)"
$ws.Cells.Item(117, 2).Value = "Wordt aan gewerkt"
$ws.Cells.Item(118, 1).Value = "Omring (Stichting)"
$ws.Cells.Item(118, 2).Value = "Actueel en vastgesteld"
$ws.Cells.Item(119, 1).Value = "Oosterlengte (Stichting)"
$ws.Cells.Item(119, 2).Value = "Wordt aan gewerkt"
$ws.Cells.Item(120, 1).Value = "Opbouw (Stichting) incl. Prinsenstichting"
$ws.Cells.Item(120, 2).Value = "Wordt aan gewerkt"
$ws.Cells.Item(121, 1).Value = "ORO (Stichting)"
$ws.Cells.Item(121, 2).Value = "Actueel en vastgesteld"
$ws.Cells.Item(122, 1).Value = "Ouderenzorg Oudewater, De Wulverhorst"
$ws.Cells.Item(122, 2).Value = "Actueel en vastgesteld"
$ws.Cells.Item(123, 1).Value = "Pameijer (Stichting)"
$ws.Cells.Item(123, 2).Value = "Actueel en vastgesteld"
$ws.Cells.Item(124, 1).Value = "Pantein (Stichting)"
$ws.Cells.Item(124, 2).Value = "Actueel en vastgesteld"
$ws.Cells.Item(125, 1).Value = "Park Zuiderhout (Stichting)"
$ws.Cells.Item(125, 2).Value = "Actueel en vastgesteld"
$ws.Cells.Item(126, 1).Value = "Parnassia Groep B.V. Incl. Parnassia haaglanden en noord holland, Antes, Brijder, Youz, etc."
$ws.Cells.Item(126, 2).Value = "Wordt aan gewerkt"
$ws.Cells.Item(127, 1).Value = "Pergamijn (Stichting)"
$ws.Cells.Item(127, 2).Value = "Wordt aan gewerkt"
$ws.Cells.Item(128, 1).Value = "Philadelphia Zorg (Stichting)"
$ws.Cells.Item(128, 2).Value = "Wordt aan gewerkt"
$ws.Cells.Item(129, 1).Value = "Pieter Raat Stichting"
$ws.Cells.Item(129, 2).Value = "Actueel en vastgesteld"
$ws.Cells.Item(130, 1).Value = "Pieter van Foreest (Stichting)"
$ws.Cells.Item(130, 2).Value = "Actueel en vastgesteld"
$ws.Cells.Item(131, 1).Value = "Pleyade (Stichting)"
$ws.Cells.Item(131, 2).Value = "Wordt aan gewerkt"
$ws.Cells.Item(132, 1).Value = "Pluryn Hoenderloo Groep (Stichting)"
$ws.Cells.Item(132, 2).Value = "Actueel en vastgesteld"
$ws.Cells.Item(133, 1).Value = "Prisma (Stichting)"
$ws.Cells.Item(133, 2).Value = "Actueel en vastgesteld"
$ws.Cells.Item(134, 1).Value = "Pro Persona"
$ws.Cells.Item(134, 2).Value = "Actueel en vastgesteld"
$ws.Cells.Item(135, 1).Value = "Pro Senectute (Stichting)"
$ws.Cells.Item(135, 2).Value = "Wordt aan gewerkt"
$ws.Cells.Item(136, 1).Value = "Profila Zorggroep (Stichting)"
$ws.Cells.Item(136, 2).Value = "Actueel en vastgesteld"
$ws.Cells.Item(137, 1).Value = "Proteion Groep (Stichting)"
$ws.Cells.Item(137, 2).Value = "Actueel en vastgesteld"
$ws.Cells.Item(138, 1).Value = "PSW (Stichting)"
$ws.Cells.Item(138, 2).Value = "Wordt aan gewerkt"
$ws.Cells.Item(139, 1).Value = "QuaRijn (Stichting)"
$ws.Cells.Item(139, 2).Value = "Wordt aan gewerkt"
$ws.Cells.Item(140, 1).Value = "R.K. Zorgcentrum Roomburgh (Stichting)"
$ws.Cells.Item(140, 2).Value = "Actueel en vastgesteld"
$ws.Cells.Item(141, 1).Value = "Raffy-Leystroom"
$ws.Cells.Item(141, 2).Value = "Wordt aan gewerkt"
$ws.Cells.Item(142, 1).Value = "RaphaelStichting"
$ws.Cells.Item(142, 2).Value = "Wordt aan gewerkt"
$ws.Cells.Item(143, 1).Value = "Residentie Molenwijck"
$ws.Cells.Item(143, 2).Value = "Actueel en vastgesteld"
$ws.Cells.Item(144, 1).Value = "Respect Zorggroep (Stichting)"
$ws.Cells.Item(144, 2).Value = "Wordt aan gewerkt"
$ws.Cells.Item(145, 1).Value = "RIBW Arnhem & Veluwevallei"
$ws.Cells.Item(145, 2).Value = "Actueel en vastgesteld"
$ws.Cells.Item(146, 1).Value = "Rijnhoven (Stichting)"
$ws.Cells.Item(146, 2).Value = "Actueel en vastgesteld"
$ws.Cells.Item(147, 1).Value = "RijnWaal Zorggroep (Stichting)"
$ws.Cells.Item(147, 2).Value = "Actueel en vastgesteld"
$ws.Cells.Item(148, 1).Value = "Rivas Zorggroep (Stichting)"
$ws.Cells.Item(148, 2).Value = "Actueel en vastgesteld"
$ws.Cells.Item(149, 1).Value = "Rivierduinen"
$ws.Cells.Item(149, 2).Value = "Wordt aan gewerkt"
$ws.Cells.Item(150, 1).Value = "Riwis Zorg & Welzijn"
$ws.Cells.Item(150, 2).Value = "Actueel en vastgesteld"
$ws.Cells.Item(151, 1).Value = "Robert Coppes Stichting"
$ws.Cells.Item(151, 2).Value = "Actueel en vastgesteld"
$ws.Cells.Item(152, 1).Value = "RST Zorgverleners, RST Zorgverleners Waardenland, RST Zorgverleners (Zwolle) (Stichting)"
$ws.Cells.Item(152, 2).Value = "Actueel en vastgesteld"
$ws.Cells.Item(153, 1).Value = "S&L Zorg (Stichting)"
$ws.Cells.Item(153, 2).Value = "Wordt aan gewerkt"
$ws.Cells.Item(154, 1).Value = "Salem Verpleeghuis (Stichting)"
$ws.Cells.Item(154, 2).Value = "Actueel en vastgesteld"
$ws.Cells.Item(155, 1).Value = "Samen Zorgen (Stichting) (ssz)"
$ws.Cells.Item(155, 2).Value = "Wordt aan gewerkt"
$ws.Cells.Item(156, 1).Value = "Santé Partners (=STMR+Vitras)"
$ws.Cells.Item(156, 2).Value = "Actueel en vastgesteld"
$ws.Cells.Item(157, 1).Value = "Schärwachter B.V."
$ws.Cells.Item(157, 2).Value = "Actueel en vastgesteld"
$ws.Cells.Item(158, 1).Value = "SEIN Stichting Epilepsie Instellingen Nederland"
$ws.Cells.Item(158, 2).Value = "Wordt aan gewerkt"
$ws.Cells.Item(159, 1).Value = "Sensire (Stichting)"
$ws.Cells.Item(159, 2).Value = "Actueel en vastgesteld"
$ws.Cells.Item(160, 1).Value = "Sevagram (Stichting)"
$ws.Cells.Item(160, 2).Value = "Actueel en vastgesteld"
$ws.Cells.Item(161, 1).Value = "Severinusstichting"
$ws.Cells.Item(161, 2).Value = "Wordt aan gewerkt"
$ws.Cells.Item(162, 1).Value = "SGL (Stichting)"
$ws.Cells.Item(162, 2).Value = "Actueel en vastgesteld"
$ws.Cells.Item(163, 1).Value = "sHeerenLoo Zorggroep (Stichting)"
$ws.Cells.Item(163, 2).Value = "Actueel en vastgesteld"
$ws.Cells.Item(164, 1).Value = "Sint Anna Boxmeer (Stichting)"
$ws.Cells.Item(164, 2).Value = "Actueel en vastgesteld"
$ws.Cells.Item(165, 1).Value = "Sint Jacob (Stichting)"
$ws.Cells.Item(165, 2).Value = "Wordt aan gewerkt"
$ws.Cells.Item(166, 1).Value = "Sint Jozef Wonen en Zorg (R.K. Stichting)"
$ws.Cells.Item(166, 2).Value = "Actueel en vastgesteld"
$ws.Cells.Item(167, 1).Value = "Siza (Stichting)"
$ws.Cells.Item(167, 2).Value = "Wordt aan gewerkt"
$ws.Cells.Item(168, 1).Value = "Solis (Stichting)"
$ws.Cells.Item(168, 2).Value = "Wordt aan gewerkt"
$ws.Cells.Item(169, 1).Value = "SOVAK (Stichting)"
$ws.Cells.Item(169, 2).Value = "Actueel en vastgesteld"
$ws.Cells.Item(170, 1).Value = "Sterk Huis (Stichting) (West Brabant was voorheen Juzt)"
$ws.Cells.Item(170, 2).Value = "Wordt aan gewerkt"
$ws.Cells.Item(171, 1).Value = "Surplus (Stichting en surplus zorg)"
$ws.Cells.Item(171, 2).Value = "Actueel en vastgesteld"
$ws.Cells.Item(172, 1).Value = "SVRZ (Stichting Voor Regionale Zorgverlening)"
$ws.Cells.Item(172, 2).Value = "Wordt aan gewerkt"
$ws.Cells.Item(173, 1).Value = "Swinhove Groep (Stichting)"
$ws.Cells.Item(173, 2).Value = "Wordt aan gewerkt"
$ws.Cells.Item(174, 1).Value = "Syndion (Stichting)"
$ws.Cells.Item(174, 2).Value = "Actueel en vastgesteld"
$ws.Cells.Item(175, 1).Value = "Tactus Verslavingszorg (Stichting)"
$ws.Cells.Item(175, 2).Value = "Wordt aan gewerkt"
$ws.Cells.Item(176, 1).Value = "tanteLouise (Stichting)"
$ws.Cells.Item(176, 2).Value = "Actueel en vastgesteld"
$ws.Cells.Item(177, 1).Value = "Teamzorg B.V."
$ws.Cells.Item(177, 2).Value = "Wordt aan gewerkt"
$ws.Cells.Item(178, 1).Value = "Teamzorg B.V."
$ws.Cells.Item(178, 2).Value = "Wordt aan gewerkt"
$ws.Cells.Item(179, 1).Value = "Terwille verslavingszorg (Stichting)"
$ws.Cells.Item(179, 2).Value = "Wordt aan gewerkt"
$ws.Cells.Item(180, 1).Value = "Thebe (Zorggroep west en midden Brabant, incl. Ruitersbos)"
$ws.Cells.Item(180, 2).Value = "Wordt aan gewerkt"
$ws.Cells.Item(181, 1).Value = "Thuis met Zorg Zaanstreek B.V."
$ws.Cells.Item(181, 2).Value = "Actueel en vastgesteld"
$ws.Cells.Item(182, 1).Value = "Topaz (Stichting)"
$ws.Cells.Item(182, 2).Value = "Wordt aan gewerkt"
$ws.Cells.Item(183, 1).Value = "Trajectum (Stichting)"
$ws.Cells.Item(183, 2).Value = "Actueel en vastgesteld"
$ws.Cells.Item(184, 1).Value = "Valkenhof (Stichting)"
$ws.Cells.Item(184, 2).Value = "Wordt aan gewerkt"
$ws.Cells.Item(185, 1).Value = "Van Neynselstichting (Stichting)"
$ws.Cells.Item(185, 2).Value = "Wordt aan gewerkt"
$ws.Cells.Item(186, 1).Value = "Vanboeijen"
$ws.Cells.Item(186, 2).Value = "Actueel en vastgesteld"
$ws.Cells.Item(187, 1).Value = "Verpleeghuis Bergweide (Stichting)"
$ws.Cells.Item(187, 2).Value = "Wordt aan gewerkt"
$ws.Cells.Item(188, 1).Value = "Viersprong (Netherlands institute for personality disorders)"
$ws.Cells.Item(188, 2).Value = "Wordt aan gewerkt"
$ws.Cells.Item(189, 1).Value = "Vilente (Stichting)"
$ws.Cells.Item(189, 2).Value = "Actueel en vastgesteld"
$ws.Cells.Item(190, 1).Value = "Vincent van Gogh (Stichting)"
$ws.Cells.Item(190, 2).Value = "Wordt aan gewerkt"
$ws.Cells.Item(191, 1).Value = "Vitalis"
$ws.Cells.Item(191, 2).Value = "Wordt aan gewerkt"
$ws.Cells.Item(192, 1).Value = "ViVa! Zorggroep (Stichting)"
$ws.Cells.Item(192, 2).Value = "Wordt aan gewerkt"
$ws.Cells.Item(193, 1).Value = "Vivent (Stichting)"
$ws.Cells.Item(193, 2).Value = "Wordt aan gewerkt"
$ws.Cells.Item(194, 1).Value = "Waardeburgh (Stichting)"
$ws.Cells.Item(194, 2).Value = "Wordt aan gewerkt"
$ws.Cells.Item(195, 1).Value = "Warande (Stichting)"
$ws.Cells.Item(195, 2).Value = "Actueel en vastgesteld"
$ws.Cells.Item(196, 1).Value = "Werkt voor Ouderen (Stichting) (WVO Zorg)"
$ws.Cells.Item(196, 2).Value = "Wordt aan gewerkt"
$ws.Cells.Item(197, 1).Value = "Wever (Stichting)"
$ws.Cells.Item(197, 2).Value = "Actueel en vastgesteld"
$ws.Cells.Item(198, 1).Value = "WIJdezorg (Stichting)"
$ws.Cells.Item(198, 2).Value = "Actueel en vastgesteld"
$ws.Cells.Item(199, 1).Value = "Wonen en Zorg Purmerend (Stichting) (SWZP)"
$ws.Cells.Item(199, 2).Value = "Wordt aan gewerkt"
$ws.Cells.Item(200, 1).Value = "Woon en zorgcentrum de Merwelanden, stichting"
$ws.Cells.Item(200, 2).Value = "Actueel en vastgesteld"
$ws.Cells.Item(201, 1).Value = "Woon- en Zorgcentrum Humanitas (Stichting)"
$ws.Cells.Item(201, 2).Value = "Actueel en vastgesteld"
$ws.Cells.Item(202, 1).Value = "Woongemeenschap voor Ouderen Heemzicht (Stichting)"
$ws.Cells.Item(202, 2).Value = "Actueel en vastgesteld"
$ws.Cells.Item(203, 1).Value = "Woonzorg Samen (Stichting)"
$ws.Cells.Item(203, 2).Value = "Actueel en vastgesteld"
$ws.Cells.Item(204, 1).Value = "Woonzorgcentrum De Zeeg (Stichting)"
$ws.Cells.Item(204, 2).Value = "Actueel en vastgesteld"
$ws.Cells.Item(205, 1).Value = "WZC de Westerkim (Stichting)"
$ws.Cells.Item(205, 2).Value = "Actueel en vastgesteld"
$ws.Cells.Item(206, 1).Value = "Youke"
$ws.Cells.Item(206, 2).Value = "Wordt aan gewerkt"
$ws.Cells.Item(207, 1).Value = "Yulius (Stichting)"
$ws.Cells.Item(207, 2).Value = "Actueel en vastgesteld"
$ws.Cells.Item(208, 1).Value = "Zellingen (Stichting Zorgbeheer De ...)"
$ws.Cells.Item(208, 2).Value = "Wordt aan gewerkt"
$ws.Cells.Item(209, 1).Value = "ZGR (Zorggroep Raalte (Stichting))"
$ws.Cells.Item(209, 2).Value = "Actueel en vastgesteld"
$ws.Cells.Item(210, 1).Value = "Zonnehuisgroep Noord (Stichting)"
$ws.Cells.Item(210, 2).Value = "Wordt aan gewerkt"
$ws.Cells.Item(211, 1).Value = "Zonnehuisgroep Vlaardingen (Stichting)"
$ws.Cells.Item(211, 2).Value = "Actueel en vastgesteld"
$ws.Cells.Item(212, 1).Value = "Zorgaccent (Stichting)"
$ws.Cells.Item(212, 2).Value = "Wordt aan gewerkt"
$ws.Cells.Item(213, 1).Value = "Zorgbalans"
$ws.Cells.Item(213, 2).Value = "Wordt aan gewerkt"
$ws.Cells.Item(214, 1).Value = "Zorgboog (Stichting)"
$ws.Cells.Item(214, 2).Value = "Actueel en vastgesteld"
$ws.Cells.Item(215, 1).Value = "Zorgcentra Rivierenland (Stichting)"
$ws.Cells.Item(215, 2).Value = "Wordt aan gewerkt"
$ws.Cells.Item(216, 1).Value = "Zorgcentrum Beek en Bos (Stichting)"
$ws.Cells.Item(216, 2).Value = "Actueel en vastgesteld"
$ws.Cells.Item(217, 1).Value = "Zorgcentrum het Bildt (Beukelaar) (Stichting)"
$ws.Cells.Item(217, 2).Value = "Wordt aan gewerkt"
$ws.Cells.Item(218, 1).Value = "Zorgcentrum 't Anker (Protestants Christelijke Stichting)"
$ws.Cells.Item(218, 2).Value = "Actueel en vastgesteld"
$ws.Cells.Item(219, 1).Value = "Zorgcentrum 't Slot (Stichting)"
$ws.Cells.Item(219, 2).Value = "Actueel en vastgesteld"
$ws.Cells.Item(220, 1).Value = "Zorgfederatie Oldenzaal (Stichting)"
$ws.Cells.Item(220, 2).Value = "Actueel en vastgesteld"
$ws.Cells.Item(221, 1).Value = "Zorggroep Amsterdam Oost (ZGAO) (Stichting)"
$ws.Cells.Item(221, 2).Value = "Actueel en vastgesteld"
$ws.Cells.Item(222, 1).Value = "Zorggroep Apeldoorn en omstreken (Stichting)"
$ws.Cells.Item(222, 2).Value = "Wordt aan gewerkt"
$ws.Cells.Item(223, 1).Value = "Zorggroep Elde Maasduinen (Maasduinen Zorg => gefuseerd uit GD HvB gestapt)"
$ws.Cells.Item(223, 2).Value = "Actueel en vastgesteld"
$ws.Cells.Item(224, 1).Value = "Zorggroep Ena (Stichting)"
$ws.Cells.Item(224, 2).Value = "Wordt aan gewerkt"
$ws.Cells.Item(225, 1).Value = "Zorggroep Groningen (Stichting)"
$ws.Cells.Item(225, 2).Value = "Actueel en vastgesteld"
$ws.Cells.Item(226, 1).Value = "Zorggroep Sint Maarten (Stichting)"
$ws.Cells.Item(226, 2).Value = "Wordt aan gewerkt"
$ws.Cells.Item(227, 1).Value = "Zorggroep Sirjon"
$ws.Cells.Item(227, 2).Value = "Actueel en vastgesteld"
$ws.Cells.Item(228, 1).Value = "Zorggroep Tellus (Stichting)"
$ws.Cells.Item(228, 2).Value = "Wordt aan gewerkt"
$ws.Cells.Item(229, 1).Value = "Zorggroep Ter Weel (Stichting)"
$ws.Cells.Item(229, 2).Value = "Actueel en vastgesteld"
$ws.Cells.Item(230, 1).Value = "Zorggroep Triade B.V. (incl. Vitree)"
$ws.Cells.Item(230, 2).Value = "Actueel en vastgesteld"
$ws.Cells.Item(231, 1).Value = "Zorgpartners Midden-Holland (Stichting)"
$ws.Cells.Item(231, 2).Value = "Actueel en vastgesteld"
$ws.Cells.Item(232, 1).Value = "ZorgSaam Zorggroep Zeeuws-Vlaanderen (Stichting)"
$ws.Cells.Item(232, 2).Value = "Actueel en vastgesteld"
$ws.Cells.Item(233, 1).Value = "ZorgSpectrum (Stichting)"
$ws.Cells.Item(233, 2).Value = "Actueel en vastgesteld"
$ws.Cells.Item(234, 1).Value = "Zorgspectrum Het Zand"
$ws.Cells.Item(234, 2).Value = "Wordt aan gewerkt"
$ws.Cells.Item(235, 1).Value = "Zozijn Beheer (Stichting)"
$ws.Cells.Item(235, 2).Value = "Wordt aan gewerkt"
$ws.Cells.Item(236, 1).Value = "ZuidOostZorg (Stichting)"
$ws.Cells.Item(236, 2).Value = "Wordt aan gewerkt"
$ws.Cells.Item(237, 1).Value = "Zuyderland Zorg (Stichting)"
$ws.Cells.Item(237, 2).Value = "Actueel en vastgesteld"
$ws.Cells.Item(238, 1).Value = "ZZG Zorggroep (Stichting)"
$ws.Cells.Item(238, 2).Value = "Wordt aan gewerkt"

# Remove the now-unused trailing rows (old sheet had data through row 263)
$ws.Range("A239:B263").Clear()

# Reset selection to A1 (matches the refreshed/default view after the data reload)
$ws.Range("A1").Select()

